# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns,
# and restore original row order for two coin pairs that had been
# swapped (rows 24/25 and 47/48), per the Dec 25 2023 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '43.877.39'
$ws.Cells.Item(2, 5).Value = '  +0.19%  '
# Row 3
$ws.Cells.Item(3, 4).Value = '2.302.35'
$ws.Cells.Item(3, 5).Value = '  +0.42%  '
# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.45%  '
# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '115.00'
$ws.Cells.Item(5, 5).Value = '  +1.35%  '
# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '266.50'
$ws.Cells.Item(6, 5).Value = '  -1.13%  '
# Row 7
$ws.Cells.Item(7, 5).Value = '  +3.41%  '
# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.01%  '
# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.617'
$ws.Cells.Item(9, 5).Value = '  -0.33%  '
# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '47.53'
$ws.Cells.Item(10, 5).Value = '  -1.36%  '
# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0943'
$ws.Cells.Item(11, 5).Value = '  -0.93%  '
# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '9.29'
$ws.Cells.Item(12, 5).Value = '  +2.91%  '
# Row 13
$ws.Cells.Item(13, 5).Value = '  +1.56%  '
# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '15.48'
$ws.Cells.Item(14, 5).Value = '  -2.09%  '
# Row 15
$ws.Cells.Item(15, 4).Value = '2.641.42'
$ws.Cells.Item(15, 5).Value = '  +0.18%  '
# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.877'
$ws.Cells.Item(16, 5).Value = '  +3.18%  '
# Row 17
$ws.Cells.Item(17, 4).Value = '2.293.80'
$ws.Cells.Item(17, 5).Value = '  +0.30%  '
# Row 18
$ws.Cells.Item(18, 4).Value = '43.677.66'
$ws.Cells.Item(18, 5).Value = '  +0.02%  '
# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.0000111'
$ws.Cells.Item(19, 5).Value = '  +0.76%  '
# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '6.81'
$ws.Cells.Item(20, 5).Value = '  +1.30%  '
# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '72.64'
$ws.Cells.Item(21, 5).Value = '  +0.47%  '
# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.45'
$ws.Cells.Item(22, 5).Value = '  -0.16%  '
# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '237.34'
$ws.Cells.Item(23, 5).Value = '  +2.03%  '
# Row 24
$ws.Cells.Item(24, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.53'
$ws.Cells.Item(24, 5).Value = '  -2.39%  '
# Row 25
$ws.Cells.Item(25, 2).Value = 'PancakeSwap'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.89'
$ws.Cells.Item(25, 5).Value = '  +3.15%  '
# Row 26
$ws.Cells.Item(26, 5).Value = '  +1.84%  '
# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '11.63'
$ws.Cells.Item(27, 5).Value = '  -0.34%  '
# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '41.48'
$ws.Cells.Item(28, 5).Value = '  -0.89%  '
# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '3.39'
$ws.Cells.Item(29, 5).Value = '  -0.43%  '
# Row 30
$ws.Cells.Item(30, 5).Value = '  -0.98%  '
# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '174.13'
$ws.Cells.Item(31, 5).Value = '  -0.66%  '
# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '21.90'
$ws.Cells.Item(32, 5).Value = '  +1.76%  '
# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0911'
$ws.Cells.Item(33, 5).Value = '  -1.27%  '
# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.70'
$ws.Cells.Item(34, 5).Value = '  +0.98%  '
# Row 35
$ws.Cells.Item(35, 5).Value = '  +2.95%  '
# Row 36
$ws.Cells.Item(36, 5).Value = '  +5.39%  '
# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.74'
$ws.Cells.Item(37, 5).Value = '  +1.56%  '
# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.91'
$ws.Cells.Item(38, 5).Value = '  +1.63%  '
# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.106'
$ws.Cells.Item(39, 5).Value = '  -1.50%  '
# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '2.59'
$ws.Cells.Item(40, 5).Value = '  +8.17%  '
# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '14.44'
$ws.Cells.Item(41, 5).Value = '  +4.89%  '
# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '75.15'
$ws.Cells.Item(42, 5).Value = '  +2.72%  '
# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.238'
$ws.Cells.Item(43, 5).Value = '  -1.65%  '
# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '6.07'
$ws.Cells.Item(44, 5).Value = '  -4.69%  '
# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.17%  '
# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '1.39'
$ws.Cells.Item(46, 5).Value = '  +0.16%  '
# Row 47
$ws.Cells.Item(47, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '1.29'
$ws.Cells.Item(47, 5).Value = '  +5.02%  '
# Row 48
$ws.Cells.Item(48, 2).Value = 'ordi'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '74.15'
$ws.Cells.Item(48, 5).Value = '  +36.93%  '
# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '8.67'
$ws.Cells.Item(49, 5).Value = '  -1.49%  '
# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.84%  '
# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '101.05'
$ws.Cells.Item(51, 5).Value = '  -1.76%  '
